# Updated symbol list on Fri Jan 13 11:44:22 UTC 2023 with GitHub Actions
# Refresh the "Price" (D) and "Volume(1h)" (E) columns for the coin rows that
# changed in this run. Values are stored as plain text (not numbers) in the
# original sheet, so each new value is written with a leading apostrophe to
# force text interpretation, then the cell style is reapplied as "Normal" so
# no stray number-format/quote-prefix styling is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "287.44"
    "E2" = "1.56%"
    "D3" = "29.50"
    "E3" = "4.08%"
    "D4" = "5.077"
    "E4" = "1.25%"
    "D5" = "0.06751"
    "E5" = "3.64%"
    "D6" = "7.363"
    "E6" = "2.05%"
    "D7" = "1.403"
    "E7" = "1.68%"
    "D8" = "0.9160"
    "E8" = "-0.28%"
    "D9" = "0.1599"
    "E9" = "3.97%"
    "D10" = "0.06941"
    "E10" = "4.54%"
    "D11" = "0.07697"
    "E11" = "1.77%"
    "D12" = "0.02923"
    "E12" = "5.05%"
    "D13" = "0.08991"
    "E13" = "0.04%"
    "D14" = "0.001572"
    "E14" = "-0.69%"
    "D15" = "0.04490"
    "E15" = "1.41%"
    "D16" = "0.0006499"
    "E16" = "1.81%"
    "D17" = "0.006163"
    "E17" = "-0.13%"
    "D18" = "3.451"
    "E18" = "0.12%"
    "D19" = "3.446"
    "E19" = "2.82%"
    "E20" = "-0.31%"
    "E21" = "0.51%"
    "D22" = "0.1307"
    "E22" = "1.97%"
    "D23" = "4.103"
    "E23" = "2.91%"
    "D24" = "0.1584"
    "E24" = "2.66%"
    "D25" = "0.001197"
    "E25" = "1.22%"
    "E26" = "-6.88%"
    "D27" = "0.0001201"
    "E27" = "0.11%"
    "D28" = "0.0001619"
    "E28" = "0.07%"
    "D40" = "0.04259"
    "E40" = "3.30%"
    "D41" = "0.006845"
    "E41" = "2.24%"
    "D42" = "0.1241"
    "E42" = "1.06%"
    "E43" = "4.82%"
    "D44" = "0.01295"
    "E44" = "7.33%"
    "D45" = "0.00005742"
    "E45" = "1.29%"
    "D46" = "1.965"
    "E46" = "-0.06%"
    "D47" = "0.01508"
    "E47" = "-18.53%"
}

foreach ($cellRef in $updates.Keys) {
    $newValue = $updates[$cellRef]
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $newValue
    $cell.Style = "Normal"
}
